$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at U (everything from the old "YD_Verify_Address" column
# onward shifts one column to the right: U->V, V->W, W->X, X->Y, ...)
$ws.Columns("U").Insert()

# New column header + value: YD_Home_Phone (U). Copy formats from the existing
# "text-like" columns (E1/E2 = YE_Estimate) so the cell styles match exactly.
$ws.Range("E1").Copy()
$ws.Range("U1").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("U2").PasteSpecial(-4122)
$ws.Range("U1").Value = "YD_Home_Phone"
$ws.Range("U2").Value = "02392123456"

# Fix YD_Verify_Address value (now shifted to column V)
$ws.Range("V2").Value = "26 Broadsands Drive,Gosport,Hampshire,PO12 2SD"

# Update YD_Email value (now shifted to column X) - keep hyperlink target as-is
$ws.Range("X2").Value = "Test_xxxxxx@xdxdxdxd.com"

# Update YD_Mobile value (now shifted to column Y) - store as text w/ leading zero
$ws.Range("E1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("Y2").PasteSpecial(-4122)
$ws.Range("Y1").Value = "YD_Mobile"
$ws.Range("Y2").Value = "07788123456"

$excel.CutCopyMode = 0

# Column width adjustments
$ws.Columns("M").ColumnWidth = 22.28515625
$ws.Columns("U").ColumnWidth = 18
$ws.Columns("V").ColumnWidth = 47.140625
$ws.Columns("Y").ColumnWidth = 13.85546875

# Fix the hyperlink anchor: it stayed on the now-stale W2 cell after the column
# insert shifted the e-mail value into X2.
$ws.Hyperlinks.Item(1).Range = $ws.Range("X2")

# Sheet view adjustments (scroll + selection)
$ws.Application.ActiveWindow.ScrollColumn = 17
$ws.Range("V2").Select()
